# Generate Report for Handback
# Applies the localization-status handback report changes:
#  - Update status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - Populate "Latest Target File" (col I) with a hyperlinked handoff file name
#  - Populate "Latest Handback File" (col J) with the locale-specific handback xliff name
#  - Populate "Latest Handback DateTime" (col K) with real handback timestamps
#  - Widen a few columns that now hold longer content

$wb = $excel.ActiveWorkbook

$zhUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e51c3425b5893d523a7283cb5d110fecc0f08f4/e2e/2b3a3bc3-03f3-4cea-bd00-5fc13465dadd.md"
$zhUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e51c3425b5893d523a7283cb5d110fecc0f08f4/e2e/d397461b-a4a6-4e6c-a039-b92e0a569f55.md"

# ---------------------------------------------------------------------------
# 1) Overview sheet: status text + column widths
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status text reused across the table
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Real handback datetime (this sheet keeps the zh-cn handback timestamp)
$wsZh.Range("K2").Value = "2016-08-31 01:06:01"
$wsZh.Range("K3").Value = "2016-08-31 01:06:01"

# Latest Target File / Latest Handback File
$wsZh.Range("J2").Value = "2b3a3bc3-03f3-4cea-bd00-5fc13465dadd.127fd7d1f1ebc4f2a8d1bd7c4469c614c3d1c54d.zh-cn.xlf"
$wsZh.Range("J3").Value = "d397461b-a4a6-4e6c-a039-b92e0a569f55.ebfefda6e4b5ec485da5b491da824fd91bf5d7a0.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhUrl1, "", "", "2b3a3bc3-03f3-4cea-bd00-5fc13465dadd.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhUrl2, "", "", "d397461b-a4a6-4e6c-a039-b92e0a569f55.md") | Out-Null

$wsZh.Range("C1").ColumnWidth = 29.17
$wsZh.Range("I1").ColumnWidth = 39.17
$wsZh.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Real handback datetime (de-de handback happened a bit later than zh-cn)
$wsDe.Range("K2").Value = "2016-08-31 01:06:14"
$wsDe.Range("K3").Value = "2016-08-31 01:06:14"

$wsDe.Range("J2").Value = "2b3a3bc3-03f3-4cea-bd00-5fc13465dadd.127fd7d1f1ebc4f2a8d1bd7c4469c614c3d1c54d.de-de.xlf"
$wsDe.Range("J3").Value = "d397461b-a4a6-4e6c-a039-b92e0a569f55.ebfefda6e4b5ec485da5b491da824fd91bf5d7a0.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $zhUrl1, "", "", "2b3a3bc3-03f3-4cea-bd00-5fc13465dadd.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $zhUrl2, "", "", "d397461b-a4a6-4e6c-a039-b92e0a569f55.md") | Out-Null

$wsDe.Range("C1").ColumnWidth = 29.17
$wsDe.Range("I1").ColumnWidth = 39.17
$wsDe.Range("J1").ColumnWidth = 39.17
